$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="4096">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Ethics in AI</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Research ethics</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>UNESCO ethics of AI</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Do no </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>harm</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Safety and security</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
        <w:t>Carbon assessments</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>Human in the loop</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Sustainable</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Fairness</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Frameworks for responsible research</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Small language models</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">BERT small language model and use these embeddings in a </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>classifier</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p/>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Huggingface</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Needle in a haystack test</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Quantization</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>QLoRa</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Quantization and LoRa</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>BERT</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Bidirectional encoder representations from transformers</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Q&amp;a</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> versions of BERT</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>CLIP</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Contrastive language image pretraining</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Infrastructure</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Runpod.io</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Model as a service</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Azure</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Bedrock</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:hyperlink r:id="rIdSTREAMLIT" w:history="1">
        <w:r>
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
          </w:rPr>
          <w:t>https://streamlit.io/</w:t>
        </w:r>
      </w:hyperlink>
    </w:p>
    <w:p>
      <w:r>
        <w:t>to make a web app</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Cursor</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Fork of VS code</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Chat models built into </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>it</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">No code </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>options</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>LMStudio</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>GPT4All</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Textgen-webui</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:hyperlink r:id="rIdOOBABOOGA" w:history="1">
        <w:r>
          <w:rPr>
            <w:rStyle w:val="Hyperlink"/>
          </w:rPr>
          <w:t>https://github.com/oobabooga/text-generation-webui</w:t>
        </w:r>
      </w:hyperlink>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Wrap-up</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
      <w:r>
        <w:t>Open and closed models</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
      <w:r>
        <w:t>Consider using a small language model</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (such as BERT)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
      <w:r>
        <w:t>Consider a quantized model</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (LORA, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>qLORA</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">RAG is </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>useful</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
      </w:pPr>
      <w:r>
        <w:t>Use temperature to control creativity</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdSTREAMLIT" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://streamlit.io/" TargetMode="External"/>
<Relationship Id="rIdOOBABOOGA" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/oobabooga/text-generation-webui" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r.InsertXML($xml)

# Apply a shared decimal numbered list ("1.", "2.", ...) to the five
# "Wrap-up" bullet paragraphs we just inserted (they already carry the
# ListParagraph style; this wires up numPr + creates numbering.xml).
$total = $d.Paragraphs.Count
$firstIdx = -1
$lastIdx = -1
for ($i = 1; $i -le $total; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Style.NameLocal -eq "List Paragraph") {
        if ($firstIdx -eq -1) { $firstIdx = $i }
        $lastIdx = $i
    }
}
if ($firstIdx -ne -1) {
    $startPos = $d.Paragraphs($firstIdx).Range.Start
    $endPos = $d.Paragraphs($lastIdx).Range.End
    $listRange = $d.Range($startPos, $endPos)
    $listRange.ListFormat.ApplyNumberDefault()
}

Write-Output "Done. Paragraphs=$($d.Paragraphs.Count) ListFirst=$firstIdx ListLast=$lastIdx"
